# "pay salary added with date"
# Rename the worksheet (this also updates the _xlnm._FilterDatabase defined name
# that refers to the sheet by name), then update the pay/balance data with the
# newly added salary payment (including its pay date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet "report-sheet" -> "employee-sheet" (definedName follows automatically)
$ws.Name = "employee-sheet"

# Row 2: update last pay date and balance for the first employee
$ws.Range("B2").Value = 44865.973287037
$ws.Range("F2").Value = 3000.5

# Row 3: new pay date uses a fresh date number format (numFmtId 14), and the
# balance now reflects the added salary payout
$ws.Range("B3").NumberFormat = "mm-dd-yy"
$ws.Range("B3").Value = 44884.986863645834
$ws.Range("F3").Value = 40000.0

# Move the active selection, matching where the editor left off
$ws.Range("E7").Select() | Out-Null
